$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.261.91'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.865.98'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''237.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '''0.4686'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").Value = '''0.2870'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.11%  '
$ws.Range("E10").Value = '  +13.30%  '
$ws.Range("D11").Value = '''0.07894'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '''97.84'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '1.869.22'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '''5.182'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '''0.6813'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").Value = '''278.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '30.265.56'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  +8.19%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''5.395'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '''0.000007347'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").Value = '2.112.92'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '''6.199'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").Value = '''168.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = '''9.294'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = '''19.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("D28").Value = '''1.942'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("D29").Value = '''1.383'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.27%  '
$ws.Range("D30").Value = '''0.09822'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("D31").Value = '''4.394'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").Value = '''4.067'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").Value = '''0.04746'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("D35").Value = '''1.143'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.80%  '
$ws.Range("D36").Value = '''0.7090'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("D37").Value = '''2.708'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("D39").Value = '''2.622'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.46%  '
$ws.Range("D40").Value = '''76.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.76%  '
$ws.Range("D41").Value = '''6.303'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").Value = '''1.963'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '
$ws.Range("D43").Value = '''0.8507'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").Value = '''0.4191'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("D45").Value = '''1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '''103.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = '''7.227'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.71%  '
$ws.Range("D48").Value = '''956.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.79%  '
$ws.Range("D49").Value = '''9.317'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").Value = '''34.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Value = '''0.05638'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.23%  '
